$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" and before "2022-Q2"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund-holding data
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3Sheet.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

$fundRows = @(
    @("007040", "新疆前海联合泳隆灵活配置混合C", "6.78", "91.50", "3.83", "0.2597", 8),
    @("001305", "九泰天富改革新动力混合A",       "1.86", "94.71", "7.74", "0.1440", 6),
    @("001844", "九泰久益灵活配置混合C",         "0.98", "93.32", "7.72", "0.0757", 6),
    @("001782", "九泰久益灵活配置混合A",         "0.53", "93.32", "7.72", "0.0409", 6),
    @("004128", "新疆前海联合泳隆灵活配置混合A", "1.02", "91.50", "3.83", "0.0391", 8),
    @("012419", "天弘国证建材指数C",             "0.63", "94.93", "4.28", "0.0270", 6),
    @("009912", "九泰天富改革新动力混合C",       "0.17", "94.71", "7.74", "0.0132", 6),
    @("007939", "华夏网购精选灵活配置混合C",     "0.41", "90.71", "1.60", "0.0066", 10),
    @("012405", "天弘国证建材指数A",             "0.13", "94.93", "4.28", "0.0056", 6),
    @("002837", "华夏网购精选灵活配置混合A",     "0.18", "90.71", "1.60", "0.0029", 10),
    @("000892", "九泰天宝灵活配置混合A",         "0.06", "89.35", "3.68", "0.0022", 9),
    @("002028", "九泰天宝灵活配置混合C",         "0.00", "89.35", "3.68", 0,        9)
)

$r = 2
foreach ($row in $fundRows) {
    $q3Sheet.Cells.Item($r, 1).Value = ($r - 2)
    $q3Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q3Sheet.Cells.Item($r, 3).Value = "'" + $row[1]
    $q3Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q3Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q3Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    if ($row[5] -eq 0) {
        $q3Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q3Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    }
    $q3Sheet.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: add a new top row for 2022-Q3 and shift
#    the existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(9).Insert()
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

$totalSheet.Cells.Item(9, 1).Value = 7
$totalSheet.Cells.Item(9, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(9, 3).Value = 51
$totalSheet.Cells.Item(9, 4).Value = 10.05

$totalSheet.Cells.Item(8, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(8, 3).Value = 29
$totalSheet.Cells.Item(8, 4).Value = 11.01

$totalSheet.Cells.Item(7, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(7, 3).Value = 38
$totalSheet.Cells.Item(7, 4).Value = 4.41

$totalSheet.Cells.Item(6, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(6, 3).Value = 81
$totalSheet.Cells.Item(6, 4).Value = 18.45

$totalSheet.Cells.Item(5, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(5, 3).Value = 34
$totalSheet.Cells.Item(5, 4).Value = 15.57

$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 17
$totalSheet.Cells.Item(4, 4).Value = 1.02

$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 30
$totalSheet.Cells.Item(3, 4).Value = 3.45

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 0.62
